$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new columns starting at column B, shifting all existing data
# (columns B onward) to the right by 9 columns.
$ws.Range("B:J").EntireColumn.Insert()

# New header values (most recent dates, newest first) for the newly
# inserted columns B1:J1.
$ws.Range("B1").Value = "Jun_16"
$ws.Range("C1").Value = "Jun_24"
$ws.Range("D1").Value = "Jun_30"
$ws.Range("E1").Value = "Jul_07"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("G1").Value = "Jul_23"
$ws.Range("H1").Value = "Aug_04"
$ws.Range("I1").Value = "Aug_25"
$ws.Range("J1").Value = "Sep_08"

# Fill "UN" for the newly inserted columns in every data row (2-33),
# matching that row's existing data extent (which itself shifted right
# by 9 columns with the insert above).
for ($r = 2; $r -le 33; $r++) {
    $addr = "B${r}:J${r}"
    $ws.Range($addr).Value = "UN"
}
